# feat(utils,command): add two functions rangeLetter and listFromColumnsStrings
# to improve gathercolums.
#
# Inserts 16 new columns (IK:IZ) before the old "email" / empty trailing
# columns on sheet1, continuing the existing 4-column repeating pattern
# (header row: Alain/Henri/Tony/Dulcinée ; data rows: OUI/NON cycle copied
# from each row's E:H block). The previous IK/IL columns (email + blank)
# shift right to JA/JB.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Insert 16 blank columns at IK:IZ - this pushes the old IK (email) and
# IL (blank numeric) columns to JA / JB respectively, and the new columns
# inherit the "s=3" style already used by the surrounding block (I:IJ).
$insertRange = $ws.Range("IK1:IZ1")
$insertRange.EntireColumn.Insert()

$firstCol = 245   # column index of "IK"
$lastRow = 9

$headerValues = @("Alain", "Henri", "Tony", "Dulcinée")

for ($c = 0; $c -lt 16; $c++) {
    $col = $firstCol + $c
    $ws.Cells.Item(1, $col).Value = $headerValues[$c % 4]
}

for ($r = 2; $r -le $lastRow; $r++) {
    # Pull this row's existing 4-value OUI/NON cycle from columns E:H (5:8)
    # so the newly inserted cells continue the same pattern as every other
    # repeated block on the row.
    $cycle = @(
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2,
        $ws.Cells.Item($r, 7).Value2,
        $ws.Cells.Item($r, 8).Value2
    )
    for ($c = 0; $c -lt 16; $c++) {
        $col = $firstCol + $c
        $ws.Cells.Item($r, $col).Value = $cycle[$c % 4]
    }
}
